$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update J74 (daily deaths for 2020-05-23) from 0 to 1
$ws.Range("J74").Value = 1

# Append the new day's row (75) right after the current last row (74).
# Copy row 73's formatting/values down to row 75 first (mirrors how the
# sheet's existing rows were built up day-by-day) so the new row inherits
# a matching look, then overwrite with the actual reported figures for
# 2020-05-24.
$ws.Range("A73:J73").Copy($ws.Range("A75:J75"))

$ws.Range("A75").Value = 43975
$ws.Range("B75").Value = 75016
$ws.Range("C75").Value = 256
$ws.Range("D75").Value = 1469
$ws.Range("E75").Value = 1
$ws.Range("F75").Value = 16
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 2
$ws.Range("I75").Value = 107
$ws.Range("J75").Value = 0

# Grow the table (and its autofilter) so the new row is included
$table = $ws.ListObjects.Item("Tabela1")
$table.Resize($ws.Range("A1:J75"))

# Leave the same single-cell selection the author ended up with
$ws.Range("J74").Select()
